$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the date NumberFormat already applied to the Fecha (D) column
$dateFormat = $ws.Cells.Item(2, 4).NumberFormat

# Rows 7-22 of the data table: Fecha(D), Volumen(J), PrecioMin(K), PrecioMax(L), PrecioProm(M), Origen(O), PrecioKg(P)
# Columns A,B,C,E,F,G,H,I,N,Q,R are identical for every data row in this sheet.
$rows = @(
    @{ Row = 7; Fecha = 44434; Volumen = 50; PMin = 15000; PMax = 15000; PProm = 15000; Origen = "Provincia de Limarí"; PKg = 600 },
    @{ Row = 8; Fecha = 44162; Volumen = 260; PMin = 7000; PMax = 8000; PProm = 7462; Origen = "Región de La Araucanía"; PKg = 298 },
    @{ Row = 9; Fecha = 44175; Volumen = 50; PMin = 8000; PMax = 8000; PProm = 8000; Origen = "Región de La Araucanía"; PKg = 320 },
    @{ Row = 10; Fecha = 44427; Volumen = 30; PMin = 15000; PMax = 15000; PProm = 15000; Origen = "Provincia de Limarí"; PKg = 600 },
    @{ Row = 11; Fecha = 44176; Volumen = 20; PMin = 11000; PMax = 11000; PProm = 11000; Origen = "Región de La Araucanía"; PKg = 440 },
    @{ Row = 12; Fecha = 44354; Volumen = 80; PMin = 16000; PMax = 16000; PProm = 16000; Origen = "Provincia de Limarí"; PKg = 640 },
    @{ Row = 13; Fecha = 44371; Volumen = 40; PMin = 15000; PMax = 15000; PProm = 15000; Origen = "Provincia de Limarí"; PKg = 600 },
    @{ Row = 14; Fecha = 44435; Volumen = 50; PMin = 15000; PMax = 15000; PProm = 15000; Origen = "Provincia de Limarí"; PKg = 600 },
    @{ Row = 15; Fecha = 44181; Volumen = 55; PMin = 14000; PMax = 14000; PProm = 14000; Origen = "Provincia de Cautín"; PKg = 560 },
    @{ Row = 16; Fecha = 44161; Volumen = 300; PMin = 7000; PMax = 7000; PProm = 7000; Origen = "Región de La Araucanía"; PKg = 280 },
    @{ Row = 17; Fecha = 44159; Volumen = 50; PMin = 8000; PMax = 8000; PProm = 8000; Origen = "Región de La Araucanía"; PKg = 320 },
    @{ Row = 18; Fecha = 44159; Volumen = 80; PMin = 8000; PMax = 8000; PProm = 8000; Origen = "Región del Maule"; PKg = 320 },
    @{ Row = 19; Fecha = 44160; Volumen = 30; PMin = 8000; PMax = 8000; PProm = 8000; Origen = "Región de La Araucanía"; PKg = 320 },
    @{ Row = 20; Fecha = 44186; Volumen = 30; PMin = 14000; PMax = 14000; PProm = 14000; Origen = "Región de La Araucanía"; PKg = 560 },
    @{ Row = 21; Fecha = 44355; Volumen = 20; PMin = 16000; PMax = 16000; PProm = 16000; Origen = "Provincia de Limarí"; PKg = 640 },
    @{ Row = 22; Fecha = 44168; Volumen = 120; PMin = 7000; PMax = 8000; PProm = 7458; Origen = "Región de La Araucanía"; PKg = 298 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 10
    $ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($row, 3).Value = "La Araucanía"
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 4).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 5).Value = 9
    $ws.Cells.Item($row, 6).Value = 100112026
    $ws.Cells.Item($row, 7).Value = "Haba"
    $ws.Cells.Item($row, 8).Value = "Sin especificar"
    $ws.Cells.Item($row, 9).Value = "Primera"
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = "$/saco 25 kilos"
    $ws.Cells.Item($row, 15).Value = $r.Origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = 25
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
